$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 29.34999999999999
$ws.Cells.Item(2, 3).Value = 29.78336143493652
$ws.Cells.Item(2, 4).Value = 0.4333614349365291
$ws.Cells.Item(2, 5).Value = 0.1878021332902476

$ws.Cells.Item(3, 3).Value = 29.6746883392334
$ws.Cells.Item(3, 4).Value = 0.3046883392333939
$ws.Cells.Item(3, 5).Value = 0.09283498406480371

$ws.Cells.Item(4, 2).Value = 29.53999999999999
$ws.Cells.Item(4, 3).Value = 29.21548271179199
$ws.Cells.Item(4, 4).Value = -0.3245172882079999
$ws.Cells.Item(4, 5).Value = 0.105311470345874

$ws.Cells.Item(5, 3).Value = 29.43997001647949
$ws.Cells.Item(5, 4).Value = -0.110029983520505
$ws.Cells.Item(5, 5).Value = 0.0121065972735226

$ws.Cells.Item(6, 3).Value = 29.63873863220215
$ws.Cells.Item(6, 4).Value = -0.1112613677978516
$ws.Cells.Item(6, 5).Value = 0.0123790919642488

$ws.Cells.Item(7, 3).Value = 29.90313529968262
$ws.Cells.Item(7, 4).Value = 0.06313529968261378
$ws.Cells.Item(7, 5).Value = 0.003986066066013452

$ws.Cells.Item(8, 3).Value = 29.88266563415527
$ws.Cells.Item(8, 4).Value = 0.07266563415527116
$ws.Cells.Item(8, 5).Value = 0.005280294387187711

$ws.Cells.Item(9, 3).Value = 29.90141868591309
$ws.Cells.Item(9, 4).Value = -0.01858131408691577
$ws.Cells.Item(9, 5).Value = 0.0003452652331966143

$ws.Cells.Item(10, 3).Value = 29.95934104919434
$ws.Cells.Item(10, 4).Value = -0.02065895080566804
$ws.Cells.Item(10, 5).Value = 0.0004267922483910122

$ws.Cells.Item(11, 2).Value = 30.03999999999999
$ws.Cells.Item(11, 3).Value = 30.1280574798584
$ws.Cells.Item(11, 4).Value = 0.0880574798584064
$ws.Cells.Item(11, 5).Value = 0.007754119759013648

$ws.Cells.Item(12, 2).Value = 30.21000000000001
$ws.Cells.Item(12, 3).Value = 30.17116355895996
$ws.Cells.Item(12, 4).Value = -0.03883644104004702
$ws.Cells.Item(12, 5).Value = 0.001508269152657048

$ws.Cells.Item(13, 3).Value = 30.2900447845459
$ws.Cells.Item(13, 4).Value = 0.07004478454589957
$ws.Cells.Item(13, 5).Value = 0.004906271842081492

$ws.Cells.Item(14, 3).Value = 30.34025382995605
$ws.Cells.Item(14, 4).Value = -0.03974617004394077
$ws.Cells.Item(14, 5).Value = 0.001579758033161854

$ws.Cells.Item(15, 3).Value = 30.54611778259277
$ws.Cells.Item(15, 4).Value = 0.1061177825927757
$ws.Cells.Item(15, 5).Value = 0.01126098378240761

$ws.Cells.Item(16, 3).Value = 30.41673851013184
$ws.Cells.Item(16, 4).Value = -0.06326148986816804
$ws.Cells.Item(16, 5).Value = 0.004002016100340328

$ws.Cells.Item(17, 3).Value = 30.42951011657715
$ws.Cells.Item(17, 4).Value = -0.2604898834228493
$ws.Cells.Item(17, 5).Value = 0.06785497936564962

$ws.Cells.Item(18, 3).Value = 30.52583885192871
$ws.Cells.Item(18, 4).Value = -0.2241611480712891
$ws.Cells.Item(18, 5).Value = 0.05024822030463838

$ws.Cells.Item(19, 3).Value = 30.64208030700684
$ws.Cells.Item(19, 4).Value = -0.2979196929931618
$ws.Cells.Item(19, 5).Value = 0.08875614347313977

$ws.Cells.Item(20, 3).Value = 30.75008964538574
$ws.Cells.Item(20, 4).Value = -0.1999103546142607
$ws.Cells.Item(20, 5).Value = 0.03996414988199944

$ws.Cells.Item(21, 3).Value = 31.1014232635498
$ws.Cells.Item(21, 4).Value = 0.08142326354980867
$ws.Cells.Item(21, 5).Value = 0.006629747847101601

$ws.Cells.Item(22, 3).Value = 31.24580955505371
$ws.Cells.Item(22, 4).Value = 0.1258095550537064
$ws.Cells.Item(22, 5).Value = 0.01582804414281158

$ws.Cells.Item(23, 3).Value = 31.31645774841309
$ws.Cells.Item(23, 4).Value = 0.0364577484130848
$ws.Cells.Item(23, 5).Value = 0.001329167419351787

$ws.Cells.Item(24, 3).Value = 31.22794914245605
$ws.Cells.Item(24, 4).Value = -0.1520508575439408
$ws.Cells.Item(24, 5).Value = 0.02311946327984777

$ws.Cells.Item(25, 3).Value = 31.41512870788574
$ws.Cells.Item(25, 4).Value = -0.1648712921142561
$ws.Cells.Item(25, 5).Value = 0.02718254296342437

$ws.Cells.Item(26, 2).Value = 31.65000000000001
$ws.Cells.Item(26, 3).Value = 31.88678550720215
$ws.Cells.Item(26, 4).Value = 0.2367855072021428
$ws.Cells.Item(26, 5).Value = 0.056067376420976

$ws.Cells.Item(27, 3).Value = 32.44120407104492
$ws.Cells.Item(27, 4).Value = 0.5612040710449264
$ws.Cells.Item(27, 5).Value = 0.3149500093573988

$ws.Cells.Item(28, 3).Value = 32.40054702758789
$ws.Cells.Item(28, 4).Value = 0.1205470275878895
$ws.Cells.Item(28, 5).Value = 0.01453158586027539

$ws.Cells.Item(29, 3).Value = 32.50535202026367
$ws.Cells.Item(29, 4).Value = 0.05535202026366903
$ws.Cells.Item(29, 5).Value = 0.003063846147269627

$ws.Cells.Item(30, 2).Value = 32.84999999999999
$ws.Cells.Item(30, 3).Value = 32.76932907104492
$ws.Cells.Item(30, 4).Value = -0.08067092895507244
$ws.Cells.Item(30, 5).Value = 0.006507798778474345

$ws.Cells.Item(31, 2).Value = 32.90000000000001
$ws.Cells.Item(31, 3).Value = 32.95795059204102
$ws.Cells.Item(31, 4).Value = 0.05795059204100994
$ws.Cells.Item(31, 5).Value = 0.003358271117903565

$ws.Cells.Item(32, 2).Value = 33.09999999999999
$ws.Cells.Item(32, 3).Value = 32.89803695678711
$ws.Cells.Item(32, 4).Value = -0.2019630432128849
$ws.Cells.Item(32, 5).Value = 0.04078907082380963

$ws.Cells.Item(33, 2).Value = 33.40000000000001
$ws.Cells.Item(33, 3).Value = 33.65871047973633
$ws.Cells.Item(33, 4).Value = 0.2587104797363224
$ws.Cells.Item(33, 5).Value = 0.0669311123253981

$ws.Cells.Item(34, 3).Value = 33.65222549438477
$ws.Cells.Item(34, 4).Value = -0.04777450561523722
$ws.Cells.Item(34, 5).Value = 0.002282403386780332

$ws.Cells.Item(35, 2).Value = 34.09999999999999
$ws.Cells.Item(35, 3).Value = 33.83626174926758
$ws.Cells.Item(35, 4).Value = -0.2637382507324162
$ws.Cells.Item(35, 5).Value = 0.06955786489939483

$ws.Cells.Item(36, 2).Value = 34.40000000000001
$ws.Cells.Item(36, 3).Value = 34.3889045715332
$ws.Cells.Item(36, 4).Value = -0.01109542846680256
$ws.Cells.Item(36, 5).Value = 0.0001231085328619326

$ws.Cells.Item(37, 2).Value = 34.90000000000001
$ws.Cells.Item(37, 3).Value = 35.00925064086914
$ws.Cells.Item(37, 4).Value = 0.1092506408691349
$ws.Cells.Item(37, 5).Value = 0.0119357025303167

$ws.Cells.Item(38, 3).Value = 35.71824264526367
$ws.Cells.Item(38, 4).Value = 0.4182426452636747
$ws.Cells.Item(38, 5).Value = 0.174926910317156

$ws.Cells.Item(39, 3).Value = 35.99924087524414
$ws.Cells.Item(39, 4).Value = 0.2992408752441378
$ws.Cells.Item(39, 5).Value = 0.08954510141687763

$ws.Cells.Item(40, 3).Value = 36.01230239868164
$ws.Cells.Item(40, 4).Value = -0.2876976013183565
$ws.Cells.Item(40, 5).Value = 0.08276990980433602

$ws.Cells.Item(41, 3).Value = 36.59271621704102
$ws.Cells.Item(41, 4).Value = -0.2072837829589815
$ws.Cells.Item(41, 5).Value = 0.04296656667778617

$ws.Cells.Item(42, 3).Value = 37.18803024291992
$ws.Cells.Item(42, 4).Value = -0.1119697570800753
$ws.Cells.Item(42, 5).Value = 0.01253722650057107

$ws.Cells.Item(43, 2).Value = 37.90000000000001
$ws.Cells.Item(43, 3).Value = 37.9639778137207
$ws.Cells.Item(43, 4).Value = 0.06397781372069744
$ws.Cells.Item(43, 5).Value = 0.004093160648480262

$ws.Cells.Item(44, 3).Value = 38.41967391967773
$ws.Cells.Item(44, 4).Value = -0.08032608032226562
$ws.Cells.Item(44, 5).Value = 0.006452279179939069

$ws.Cells.Item(45, 2).Value = 38.90000000000001
$ws.Cells.Item(45, 3).Value = 39.0008659362793
$ws.Cells.Item(45, 4).Value = 0.1008659362792912
$ws.Cells.Item(45, 5).Value = 0.01017393710149803

$ws.Cells.Item(46, 2).Value = 39.40000000000001
$ws.Cells.Item(46, 3).Value = 39.5327262878418
$ws.Cells.Item(46, 4).Value = 0.1327262878417912
$ws.Cells.Item(46, 5).Value = 0.01761626748426201

$ws.Cells.Item(47, 2).Value = 39.90000000000001
$ws.Cells.Item(47, 3).Value = 39.54935836791992
$ws.Cells.Item(47, 4).Value = -0.3506416320800838
$ws.Cells.Item(47, 5).Value = 0.1229495541477849

$ws.Cells.Item(48, 2).Value = 40.09999999999999
$ws.Cells.Item(48, 3).Value = 40.00116348266602
$ws.Cells.Item(48, 4).Value = -0.09883651733397869
$ws.Cells.Item(48, 5).Value = 0.00976865715870987

$ws.Cells.Item(49, 2).Value = 40.59999999999999
$ws.Cells.Item(49, 3).Value = 40.52347183227539
$ws.Cells.Item(49, 4).Value = -0.07652816772460369
$ws.Cells.Item(49, 5).Value = 0.005856560455285074

$ws.Cells.Item(50, 2).Value = 40.90000000000001
$ws.Cells.Item(50, 3).Value = 40.82853317260742
$ws.Cells.Item(50, 4).Value = -0.07146682739258381
$ws.Cells.Item(50, 5).Value = 0.005107507417561367

$ws.Cells.Item(51, 2).Value = 41.20000000000001
$ws.Cells.Item(51, 3).Value = 41.49383163452148
$ws.Cells.Item(51, 4).Value = 0.2938316345214744
$ws.Cells.Item(51, 5).Value = 0.08633702944556132

$ws.Cells.Item(52, 3).Value = 0.1741580963134552
$ws.Cells.Item(52, 5).Value = 2.03359539015778

$ws.Cells.Item(53, 5).Value = 0.0406719078031556
